$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs: edit only the characters that changed) ---
# "Volume 30   Number  24" -> "...25"  (A8, merged A8:B8)
$ws.Range("A8").Characters(21, 2).Text = "25"
# "Report Covering the Week  6/12/2023  Through  6/18/2023" -> 6/19/2023 ... 6/25/2023 (C9, merged C9:L9)
$ws.Range("C9").Characters(27, 9).Text = "6/19/2023"
$ws.Range("C9").Characters(47, 9).Text = "6/25/2023"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -94.117647058823
$ws.Range("F15").Value = 2
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -41.176470588235
$ws.Range("L15").Value = -37.5
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -71.428571428571
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 157
$ws.Range("J16").Value = 137
$ws.Range("K16").Value = 14.598540145985
$ws.Range("L16").Value = 41.441441441441
$ws.Range("M16").Value = -13.259668508287
$ws.Range("N16").Value = -77.507163323782
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -2.564102564102
$ws.Range("I17").Value = 211
$ws.Range("J17").Value = 162
$ws.Range("K17").Value = 30.246913580246
$ws.Range("L17").Value = 47.552447552447
$ws.Range("M17").Value = 24.852071005917
$ws.Range("N17").Value = -48.661800486618
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -8
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 147
$ws.Range("K18").Value = -26.530612244898
$ws.Range("L18").Value = -12.195121951219
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -82.352941176470
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 18
$ws.Range("I19").Value = 335
$ws.Range("J19").Value = 307
$ws.Range("K19").Value = 9.120521172638
$ws.Range("L19").Value = 51.583710407239
$ws.Range("M19").Value = 141.007194244604
$ws.Range("N19").Value = 30.350194552529
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 38.888888888888
$ws.Range("I20").Value = 96
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = -4.950495049504
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 31.506849315068
$ws.Range("N20").Value = -79.617834394904
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 180
$ws.Range("G21").Value = 167
$ws.Range("H21").Value = 7.784431137724
$ws.Range("I21").Value = 918
$ws.Range("J21").Value = 873
$ws.Range("K21").Value = 5.154639175257
$ws.Range("L21").Value = 33.236574746008
$ws.Range("M21").Value = 23.553162853297
$ws.Range("N21").Value = -63.294682127149
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 166.666666666667
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = 60
$ws.Range("I23").Value = 21
$ws.Range("K23").Value = 61.538461538461
$ws.Range("L23").Value = 31.25
$ws.Range("M23").Value = 250
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 153.333333333333
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 31.707317073170
$ws.Range("I24").Value = 491
$ws.Range("J24").Value = 472
$ws.Range("K24").Value = 4.025423728813
$ws.Range("L24").Value = 19.174757281553
$ws.Range("M24").Value = 43.988269794721
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = -3.571428571428
$ws.Range("I25").Value = 325
$ws.Range("J25").Value = 302
$ws.Range("K25").Value = 7.615894039735
$ws.Range("L25").Value = 54.761904761904
$ws.Range("M25").Value = -14.248021108179
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = -9.090909090909
$ws.Range("L26").Value = -25.925925925925
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -12.5
$ws.Range("M28").Value = -73.684210526315
$ws.Range("N28").Value = -94.845360824742
$ws.Range("M29").Value = -64.285714285714
$ws.Range("N29").Value = -94.623655913978

# --- Cells that flip from a number to a text marker ("0" or "***.*") ---
# Force text storage first (so "0" is not silently re-parsed as numeric 0), then
# restore the original General/right-aligned look by pasting formats from a
# same-style neighbour that the diff leaves untouched (style 14).
$fmtSrc14 = $ws.Range("C14")
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$fmtSrc14.Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$fmtSrc14.Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$fmtSrc14.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$fmtSrc14.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$fmtSrc14.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$fmtSrc14.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$fmtSrc14.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$fmtSrc14.Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cells that flip from a text marker to a real number ---
# Plain assignment switches the stored type to Number but keeps the old text-style
# (14); paste the matching numeric style back in from a stable neighbour.
$fmtSrc15 = $ws.Range("I14")   # style 15: #,##0 integer
$fmtSrc16 = $ws.Range("K14")   # style 16: #,##0.0 percent-change
$ws.Range("D15").Value = 2
$fmtSrc15.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$fmtSrc16.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$fmtSrc15.Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$fmtSrc16.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("D26").Value = 2
$fmtSrc15.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = -50
$fmtSrc16.Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
